$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original column A width (character units) before inserting,
# so the freshly-inserted column can be resized to match it.
$origColAWidth = $ws.Range("A1").ColumnWidth

# Insert a new column before column B. Excel shifts the old column B
# ("dbExcel" / Neo4jData file name) to C and old column C ("WebExcel" /
# WebData file name) to D, carrying their widths/content with them.
$ws.Range("B1").EntireColumn.Insert()

# Header for the newly inserted StatQuery column
$ws.Range("B1").Value = "StatQuery"

# Give the new column the same width as column A (mirrors what Excel does
# when a column is inserted immediately to the right of a formatted one).
$ws.Range("B1").ColumnWidth = $origColAWidth

# New stat-bar query text for row 2, with the same wrap-text styling as A2
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$ws.Range("B2").Value = $statQuery
$ws.Range("B2").WrapText = $true

# Update the saved selection to A2, matching the workbook's saved view state
$ws.Range("A2").Select()
